# Load Screen code Updated
#
# Applies the same edits as the authoritative OOXML diff:
#  - "Add Load" sheet: row3 "Changepond T" -> "New Day"; add a new row4
#    (a copy of row3, but for test case "Loads_TC003")
#  - "Edit Load" sheet: row3 "ChanDriver625831" -> "Karthik"; row4
#    "ORGIN WEIGHT" -> "ORIGIN WEIGHT" (typo fix); add new rows 6-8
#    (copies of rows 3-5, but for test case "Loads_TC003")
#  - Re-point the active tab / selections: "Edit Load" becomes the active
#    sheet with C7 selected; "Add Load" ends with A7 selected; "View Load"
#    ends with E5 selected (and loses tabSelected).

$wb = $excel.ActiveWorkbook

$wsAdd  = $wb.Worksheets.Item("Add Load")
$wsEdit = $wb.Worksheets.Item("Edit Load")
$wsView = $wb.Worksheets.Item("View Load")

# ---------------------------------------------------------------------
# "Add Load" sheet - row 3 edit
# ---------------------------------------------------------------------

# Row 3: rename the driver/test value from "Changepond T" to "New Day"
$wsAdd.Range("B3").Value = "New Day"

# ---------------------------------------------------------------------
# "Edit Load" sheet - row 3 / row 4 edits
# ---------------------------------------------------------------------

# Row 3: driver name value fixed
$wsEdit.Range("C3").Value = "Karthik"

# Row 4: fix "ORGIN WEIGHT" typo -> "ORIGIN WEIGHT"
$wsEdit.Range("B4").Value = "ORIGIN WEIGHT"

# ---------------------------------------------------------------------
# "Add Load" sheet - new row 4 (duplicate of row 3, for Loads_TC003)
# ---------------------------------------------------------------------

$wsAdd.Range("A4").Value = "Loads_TC003"
$wsAdd.Range("B4").Value = "New Day"
$wsAdd.Range("C4").Value = "Current Date"
$wsAdd.Range("D4").Value = "CP Shipper"
$wsAdd.Range("E4").Value = "TestContact"
$wsAdd.Range("F4").Value = "Corn"
$wsAdd.Range("G4").NumberFormat = "@"
$wsAdd.Range("G4").Value = "0.25"
$wsAdd.Range("H4").Value = "Bushels"
$wsAdd.Range("I4").Value = "Alaska"
$wsAdd.Range("J4").Value = "Roger"
$wsAdd.Range("K4").Value = "Added new load successfully"

# ---------------------------------------------------------------------
# "Edit Load" sheet - new rows 6-8 (duplicates of rows 3-5, for Loads_TC003)
# ---------------------------------------------------------------------

# New row 6 - duplicate of row 3, but for Loads_TC003
$wsEdit.Range("A6").Value = "Loads_TC003"
$wsEdit.Range("B6").Value = "Driver"
$wsEdit.Range("C6").Value = "Karthik"
$wsEdit.Range("D6").Value = "Load Edited Successfully"

# New row 7 - duplicate of row 4, but for Loads_TC003
$wsEdit.Range("A7").Value = "Loads_TC003"
$wsEdit.Range("B7").Value = "ORIGIN WEIGHT"
$wsEdit.Range("C7").Value = 5000
$wsEdit.Range("C7").NumberFormat = "@"
$wsEdit.Range("D7").Value = "Load Edited Successfully"

# New row 8 - duplicate of row 5, but for Loads_TC003
$wsEdit.Range("A8").Value = "Loads_TC003"
$wsEdit.Range("B8").Value = "DESTINATION WEIGHt"
$wsEdit.Range("C8").NumberFormat = "@"
$wsEdit.Range("C8").Value = "5000"
$wsEdit.Range("D8").Value = "Load Edited Successfully"

# ---------------------------------------------------------------------
# Selections / active sheet
# ---------------------------------------------------------------------

# "Add Load" ends up with A7 selected
[void]$wsAdd.Range("A7").Select()

# "View Load" ends up with E5 selected (and is no longer the active tab)
[void]$wsView.Range("E5").Select()

# "Edit Load" becomes the active sheet/tab with C7 selected - select this
# last so it ends up being the active tab in the saved workbook.
[void]$wsEdit.Activate()
[void]$wsEdit.Range("C7").Select()
